# Updated stencil timings for srlpi.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header labels (E1/F1/G1) ---
# The "stencil-1947ae" matrix was renamed to "stencil-4b82ec";
# its speedup column label changes accordingly, and "NumStencils"
# keeps its place in the shared-string table ahead of the new names.
$ws.Range("E1").Value = "stencil-4b82ec"
$ws.Range("F1").Value = "stencil-4b82ec speedup"
$ws.Range("G1").Value = "NumStencils"

# --- Update column E (new stencil timings) for rows 2-24 ---
$eValues = @{
  2  = 71825.5
  3  = 41537
  4  = 48000.5
  5  = 51993.8
  6  = 119228
  7  = 343055
  8  = 18202.3
  9  = 20364.3
  10 = 40897.699999999997
  11 = 95015.8
  12 = 30332
  13 = 39338.699999999997
  14 = 71267.899999999994
  15 = 202234
  16 = 71218.399999999994
  17 = 79437.5
  18 = 127481
  19 = 163021
  20 = 295181
  21 = 91989.4
  22 = 88541.2
  23 = 157507
  24 = 293535
}

foreach ($row in $eValues.Keys) {
    $ws.Cells.Item($row, 5).Value = $eValues[$row]
}

# --- Re-apply D and F column formulas as shared formulas ---
# (values recompute automatically from the updated E column)
$ws.Range("D2:D24").Formula = "=B2/C2"
$ws.Range("F2:F24").Formula = "=B2/E2"

# --- Update the active selection shown in the sheet view ---
$ws.Range("E1").Select()
